$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text number format on D/E columns so values assigned as strings
# (e.g. "0.9994", "239.54") are stored as text, matching the source workbook,
# instead of being auto-coerced to numeric cells by Excel.
$textCells = @("D2", "E2", "D3", "E3", "D4", "E4", "D5", "E5", "D6", "E6", "E7", "D8", "D9", "E9", "D10", "E10", "E11", "D12", "E12", "D13", "E13", "D14", "E14", "D15", "E15", "D16", "E16", "D17", "D18", "E18", "D19", "E19", "D20", "E20", "D21", "E21", "D22", "E22", "D23", "E23", "D24", "E24", "E25", "D26", "E26", "D27", "E27", "E28", "D29", "E29", "D30", "E30", "D31", "E31", "D32", "E32", "D33", "E33", "D34", "E34", "D35", "E35", "D36", "E36", "D37", "E37", "E38", "D39", "E39", "D40", "E40", "D41", "E41", "D42", "E42", "B43", "C43", "D43", "E43", "B44", "C44", "D44", "E44", "B45", "C45", "D45", "E45", "D46", "E46", "D47", "E47", "B48", "C48", "D48", "E48", "B49", "C49", "D49", "E49", "B50", "C50", "D50", "E50", "D51", "E51")
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range("D2").Value = "29.027.93"
$ws.Range("E2").Value = "  -1.95%  "
$ws.Range("D3").Value = "1.832.84"
$ws.Range("E3").Value = "  -1.40%  "
$ws.Range("D4").Value = "0.9994"
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "239.54"
$ws.Range("E5").Value = "  -2.26%  "
$ws.Range("D6").Value = "0.6705"
$ws.Range("E6").Value = "  -3.41%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("D8").Value = "0.07409"
$ws.Range("D9").Value = "0.2948"
$ws.Range("E9").Value = "  -3.68%  "
$ws.Range("D10").Value = "22.68"
$ws.Range("E10").Value = "  -4.52%  "
$ws.Range("E11").Value = "  -1.60%  "
$ws.Range("D12").Value = "1.830.12"
$ws.Range("E12").Value = "  -1.46%  "
$ws.Range("D13").Value = "4.999"
$ws.Range("E13").Value = "  -2.89%  "
$ws.Range("D14").Value = "0.6717"
$ws.Range("E14").Value = "  -2.99%  "
$ws.Range("D15").Value = "86.24"
$ws.Range("E15").Value = "  -5.66%  "
$ws.Range("D16").Value = "6.132"
$ws.Range("E16").Value = "  -6.59%  "
$ws.Range("D17").Value = "29.040.96"
$ws.Range("D18").Value = "0.000008214"
$ws.Range("E18").Value = "  -0.93%  "
$ws.Range("D19").Value = "227.00"
$ws.Range("E19").Value = "  -5.43%  "
$ws.Range("D20").Value = "12.41"
$ws.Range("E20").Value = "  -2.77%  "
$ws.Range("D21").Value = "0.9994"
$ws.Range("E21").Value = "  +0.01%  "
$ws.Range("D22").Value = "7.310"
$ws.Range("E22").Value = "  -3.99%  "
$ws.Range("D23").Value = "1.0000"
$ws.Range("E23").Value = "  +0.00%  "
$ws.Range("D24").Value = "160.06"
$ws.Range("E24").Value = "  +0.25%  "
$ws.Range("E25").Value = "  -4.75%  "
$ws.Range("D26").Value = "8.662"
$ws.Range("E26").Value = "  -2.99%  "
$ws.Range("D27").Value = "17.94"
$ws.Range("E27").Value = "  -1.90%  "
$ws.Range("E28").Value = "  -2.33%  "
$ws.Range("D29").Value = "4.229"
$ws.Range("E29").Value = "  -0.61%  "
$ws.Range("D30").Value = "4.110"
$ws.Range("E30").Value = "  -1.62%  "
$ws.Range("D31").Value = "1.196"
$ws.Range("E31").Value = "  -0.47%  "
$ws.Range("D32").Value = "0.05355"
$ws.Range("E32").Value = "  +3.75%  "
$ws.Range("D33").Value = "0.7492"
$ws.Range("E33").Value = "  -2.78%  "
$ws.Range("D34").Value = "1.849"
$ws.Range("E34").Value = "  -2.33%  "
$ws.Range("D35").Value = "1.122"
$ws.Range("E35").Value = "  -2.65%  "
$ws.Range("D36").Value = "2.684"
$ws.Range("E36").Value = "  -0.13%  "
$ws.Range("D37").Value = "1.291.82"
$ws.Range("E37").Value = "  -3.23%  "
$ws.Range("E38").Value = "  -3.65%  "
$ws.Range("D39").Value = "2.707"
$ws.Range("E39").Value = "  -0.79%  "
$ws.Range("D40").Value = "0.9231"
$ws.Range("E40").Value = "  -4.88%  "
$ws.Range("D41").Value = "6.060"
$ws.Range("E41").Value = "  +4.02%  "
$ws.Range("D42").Value = "0.00000000130"
$ws.Range("E42").Value = "  +3.64%  "
$ws.Range("B43").Value = "XinFinNetwork"
$ws.Range("C43").Value = "https://coinranking.com/coin/77jGXSqWJ1ofG+xinfinnetwork-xdc"
$ws.Range("D43").Value = "0.08293"
$ws.Range("E43").Value = "  +27.24%  "
$ws.Range("B44").Value = "Quant"
$ws.Range("C44").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D44").Value = "104.11"
$ws.Range("E44").Value = "  -2.45%  "
$ws.Range("B45").Value = "PaxDollar"
$ws.Range("C45").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D45").Value = "0.9992"
$ws.Range("E45").Value = "  -0.02%  "
$ws.Range("D46").Value = "1.973.99"
$ws.Range("E46").Value = "  -1.35%  "
$ws.Range("D47").Value = "0.5177"
$ws.Range("E47").Value = "  -0.77%  "
$ws.Range("B48").Value = "Aave"
$ws.Range("C48").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D48").Value = "63.50"
$ws.Range("E48").Value = "  +0.10%  "
$ws.Range("B49").Value = "RenderToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D49").Value = "1.748"
$ws.Range("E49").Value = "  -1.72%  "
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").Value = "9.318"
$ws.Range("E50").Value = "  -4.91%  "
$ws.Range("D51").Value = "0.05922"
$ws.Range("E51").Value = "  -0.14%  "
